# "Generate Report for Handback"
#
# The handback transform for the 5fc53497-... file failed (the returned
# handback file name didn't match the expected handoff file name), so the
# localization-status report is regenerated to:
#   - flip that file's Status from "Ready for handoff" to
#     "Handback transform failed" (shown both on the Overview sheet, once
#     per locale column, and on each locale's own Status column), and
#   - record the mismatch detail in the (previously empty) Error Detail
#     column for that row on each locale sheet, widening that column so the
#     message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 5fc53497-... file; E3 = zh-cn status,
# F3 = de-de status.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Each locale sheet repeats the same Status value in column C, row 3.
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Error Detail (column P) for row 3 on each locale sheet.
$zhcn.Range("P3").Value = "Handback file name: zwwgznra.s3u is different with handoff file name: 5fc53497-b23f-444b-8572-23e5740c3001.ec664498a9f266d63e2faa670e42ba9ebf5dd7b3.zh-cn."
$dede.Range("P3").Value = "Handback file name: zwwgznra.s3u is different with handoff file name: 5fc53497-b23f-444b-8572-23e5740c3001.ec664498a9f266d63e2faa670e42ba9ebf5dd7b3.de-de."

# Widen the Error Detail column (P / 16) on both locale sheets now that it
# holds real text.
$zhcn.Columns.Item(16).ColumnWidth = 39.14
$dede.Columns.Item(16).ColumnWidth = 39.14
